{"js": "// Replace the placeholder \"Diretoria de Ensino\" tokens with their final\n// values, both in the body (the salutation \"A QWREW,\") and in the page\n// header (the address block repeated below the letterhead).\n//\n// Each token is searched (case-sensitive, whole word) within the given\n// range and every match is replaced in place.\nasync function replaceAllWholeWord(range, find, replaceWith) {\n  const results = range.search(find, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const result of results.items) {\n    result.insertText(replaceWith, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1) Body: \"A QWREW,\" -> \"A QWR,\"\nawait replaceAllWholeWord(context.document.body, \"QWREW\", \"QWR\");\n\n// 2) Header (primary header of the first/only section): several\n//    placeholder tokens used throughout the address block.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst header = sections.items[0].getHeader(\"Primary\");\n\nawait replaceAllWholeWord(header, \"REW\", \"QWER\");\nawait replaceAllWholeWord(header, \"QWREW\", \"QWR\");\nawait replaceAllWholeWord(header, \"Rew\", \"Qwer\");\nawait replaceAllWholeWord(header, \"rew\", \"qwer\");\n\nawait context.sync();\n", "ps1": "# Replace the placeholder \"Diretoria de Ensino\" tokens with their final\n# values, both in the body (the salutation \"A QWREW,\") and in the page\n# header (the address block repeated below the letterhead).\n\n$d = $word.ActiveDocument\n\n# wdReplaceAll\n$wdReplaceAll = 2\n# wdFindContinue (keep searching to the end of the range instead of stopping)\n$wdFindContinue = 1\n\nfunction Replace-WholeWord($range, [string]$findText, [string]$replaceText) {\n    $range.Find.Execute(\n        $findText,    # FindText\n        $true,        # MatchCase\n        $true,        # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        $wdFindContinue, # Wrap\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        $wdReplaceAll # Replace\n    )\n}\n\n# 1) Body: \"A QWREW,\" -> \"A QWR,\"\nReplace-WholeWord $d.Content \"QWREW\" \"QWR\"\n\n# 2) Header (primary header of the first/only section): several\n#    placeholder tokens used throughout the address block.\n$header = $d.Sections.Item(1).Headers.Item(1)\n\nReplace-WholeWord $header.Range \"REW\" \"QWER\"\nReplace-WholeWord $header.Range \"QWREW\" \"QWR\"\nReplace-WholeWord $header.Range \"Rew\" \"Qwer\"\nReplace-WholeWord $header.Range \"rew\" \"qwer\"\n"}
